{"js": "// Replace the sentence about free attendance / sponsorship funding usage\n// with the updated wording about using sponsorship funding to support\n// students through scholarships.\nconst oldText =\n  \" plan to have the attendance free for students, using the sponsorship funding for covering organization expenses.\";\nconst newText =\n  \" plan to use the sponsorship funding for supporting students through scholarships.\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found in document body.\");\n}\n\n// Replace in place so the surrounding run (and its formatting: not bold,\n// 12pt) is preserved; Word itself will re-split/merge runs as needed.\nresults.items[0].insertText(newText, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Replace the sentence about free attendance / sponsorship funding usage\n# with the updated wording about using sponsorship funding to support\n# students through scholarships.\n$d = $word.ActiveDocument\n\n$oldText = \" plan to have the attendance free for students, using the sponsorship funding for covering organization expenses.\"\n$newText = \" plan to use the sponsorship funding for supporting students through scholarships.\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.MatchCase = $true\n$find.MatchWildcards = $false\n$find.Execute()\n\nif ($find.Found) {\n    $find.Parent.Text = $newText\n}\n"}
